$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 14:04"

# Row 4
$ws.Range("B4").Value = 2234963
$ws.Range("C4").Value = 492
$ws.Range("E4").Value = 1196224

# Row 7
$ws.Range("B7").Value = 368705
$ws.Range("C7").Value = 1441
$ws.Range("D7").Value = 195139
$ws.Range("E7").Value = 161286
$ws.Range("G7").Value = 18
$ws.Range("H7").Value = 12280

# Row 12
$ws.Range("A12").Value = "Chile"
$ws.Range("B12").Value = 220628
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 156232
$ws.Range("E12").Value = 60781
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 3615

# Row 13
$ws.Range("A13").Value = "Iran"
$ws.Range("B13").Value = 197647
$ws.Range("C13").Value = 2596
$ws.Range("D13").Value = 156991
$ws.Range("E13").Value = 31384
$ws.Range("G13").Value = 87
$ws.Range("H13").Value = 9272

# Row 14
$ws.Range("A14").Value = "Alemania"
$ws.Range("B14").Value = 190179
$ws.Range("D14").Value = 173600
$ws.Range("E14").Value = 7652
$ws.Range("H14").Value = 8927

# Row 27
$ws.Range("B27").Value = 56657
$ws.Range("C27").Value = 625
$ws.Range("D27").Value = 34023
$ws.Range("E27").Value = 22303
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 331

# Row 32
$ws.Range("B32").Value = 43752
$ws.Range("C32").Value = 388
$ws.Range("D32").Value = 30241
$ws.Range("E32").Value = 13213
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 298

# Row 35
$ws.Range("A35").Value = "Kuwait"
$ws.Range("B35").Value = 38074
$ws.Range("C35").Value = 541
$ws.Range("D35").Value = 29512
$ws.Range("E35").Value = 8254
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 308

# Row 36
$ws.Range("A36").Value = "Portugal"
$ws.Range("B36").Value = 37672
$ws.Range("D36").Value = 23580
$ws.Range("E36").Value = 12569
$ws.Range("H36").Value = 1523

# Row 55
$ws.Range("B55").Value = 17223
$ws.Range("C55").Value = 20
$ws.Range("D55").Value = 16101
$ws.Range("E55").Value = 434
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 688

# Row 60
$ws.Range("B60").Value = 12334
$ws.Range("C60").Value = 40
$ws.Range("D60").Value = 11242
$ws.Range("E60").Value = 492
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 600

# Row 69
$ws.Range("E69").Value = 310
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 244

# Row 72
$ws.Range("A72").Value = "Nepal"
$ws.Range("B72").Value = 7848
$ws.Range("C72").Value = 671
$ws.Range("D72").Value = 1186
$ws.Range("E72").Value = 6640
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 22

# Row 73
$ws.Range("A73").Value = "Australia"
$ws.Range("B73").Value = 7391
$ws.Range("C73").Value = 21
$ws.Range("D73").Value = 6877
$ws.Range("E73").Value = 412
$ws.Range("H73").Value = 102

# Row 104
$ws.Range("B104").Value = 1926
$ws.Range("C104").Value = 2
$ws.Range("E104").Value = 494

# Row 120
$ws.Range("B120").Value = 1403
$ws.Range("C120").Value = 25
$ws.Range("D120").Value = 463
$ws.Range("E120").Value = 927
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 13

# Row 130
$ws.Range("D130").Value = 810
$ws.Range("E130").Value = 36

# Row 137
$ws.Range("D137").Value = 486
$ws.Range("E137").Value = 255

# Row 139
$ws.Range("D139").Value = 609
$ws.Range("E139").Value = 42
$ws.Range("G139").Value = 3
$ws.Range("H139").Value = 45

# Row 141
$ws.Range("B141").Value = 663
$ws.Range("C141").Value = 1
$ws.Range("E141").Value = 44

# Row 156
$ws.Range("A156").Value = "Vietnam"
$ws.Range("B156").Value = 342
$ws.Range("C156").Value = 7
$ws.Range("E156").Value = 17
$ws.Range("H156").Value = 0

# Row 157
$ws.Range("A157").Value = "Mauricio"
$ws.Range("B157").Value = 337
$ws.Range("D157").Value = 325
$ws.Range("E157").Value = 2
$ws.Range("H157").Value = 10

# Row 158
$ws.Range("A158").Value = "Isla de Man"
$ws.Range("B158").Value = 336
$ws.Range("D158").Value = 312
$ws.Range("E158").Value = 0
$ws.Range("H158").Value = 24

# Row 162
$ws.Range("B162").Value = 221
$ws.Range("C162").Value = 19
$ws.Range("E162").Value = 109

# Row 202
$ws.Range("A202").Value = "Fiyi"

# Row 203
$ws.Range("A203").Value = "Dominica"

# Row 208
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
